$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Repayment schedule")
$ws.Activate()
$ws.Columns("N").Insert()
$ws.Columns("N").ColumnWidth = 10.7109375
$ws.Range("R7").Select()
